$wb = $excel.ActiveWorkbook

# --- "Math & Algorithms" sheet (physically sheet1.xml): populate the new
# tracked-progress row (row 30) of the second table (I1:N30 / "Table4") ---
$ws1 = $wb.Worksheets.Item("Math & Algorithms")
$ws1.Range("I30").Value = "Problem Solving(Algorithms & Data Structures)"
$ws1.Range("J30").Value = 45159
$ws1.Range("L30").Value = "2068.16/2200"
$ws1.Range("K30").Value = "90% (131.84 points to next star)"
$ws1.Range("M30").Value = 49342
$ws1.Range("N30").Formula = "=IF(ROW()>2,(`$M`$2-M30)/`$M`$2,""NA"")"

# --- "Plan" sheet: fill in newly-planned / in-progress problems ---
$ws3 = $wb.Worksheets.Item("Plan")

# "Difficult Problems" section currently being attempted
$ws3.Range("A7").Value = "Queen's Attack II"
$ws3.Range("B7").Value = "started working on it  -almost done - need extra testing & validation - current solution working fine for some cases but not for all"
$ws3.Range("B7").WrapText = $true
$ws3.Rows("7").RowHeight = 90

# Next problems to solve
$ws3.Range("A18").Value = "Fibonacci Finding (easy)"
$ws3.Range("A20").Value = "Points on a Rectangle"

$ws3.Range("D4").Value = "Hackerland Radio Transmitters`n"
$ws3.Range("D4").WrapText = $true
$ws3.Rows("4").RowHeight = 45

$ws3.Range("A19").Value = "Special Multiple"
$ws3.Range("B19").Value = "start with 9 and generate two other numbers for each number (use tree) ex.: 90 and 99 then 900&909 AND 990&999,… and so on"
$ws3.Range("B19").WrapText = $true
$ws3.Rows("19").RowHeight = 30

# Mark "Bigger is Greater" as done/solved (strike-through, like "The Time in Words")
$ws3.Range("A5").Font.Strikethrough = $true

# Column widths: comments column (B) now holds the long free-text notes,
# category column (A) just holds short problem names.
$ws3.Columns("A").ColumnWidth = 27.85546875
$ws3.Columns("B").ColumnWidth = 62.28515625
$ws3.Columns("D").ColumnWidth = 24.28515625
$ws3.Columns("E").ColumnWidth = 62.28515625

# --- Selection / active sheet bookkeeping, matching the saved view state ---
$ws1.Range("I2:N3").Select()
$ws1.Application.ActiveWindow.ScrollRow = 1

$ws2 = $wb.Worksheets.Item("Python,C++ & SQL")
$ws2.Range("I9").Select()

$ws3.Range("B19").Select()
$ws3.Activate()
